$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23, shifting existing rows 23-27 down to 24-28
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new data record
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "Macroferia Regional de Talca"
$ws.Range("C23").Value = "Maule"
$ws.Range("D23").Value = 44468
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = 100112026
$ws.Range("G23").Value = "Haba"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 9000
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = 9000
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Región de O'Higgins"
$ws.Range("P23").Value = 360
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
